$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1359.8
$ws.Range("I40").Value = 1250
$ws.Range("J40").Value = 1799
$ws.Range("K40").Value = 1250
$ws.Range("L40").Value = 1799
$ws.Range("M40").Value = -1075
$ws.Range("N40").Value = -2149
# Row 43
$ws.Range("H43").Value = 3262.8572
$ws.Range("I43").Value = 1971.4286
$ws.Range("J43").Value = 4554.2856
$ws.Range("K43").Value = 1971.4286
$ws.Range("L43").Value = 4554.2856
$ws.Range("M43").Value = -1902.4286
$ws.Range("N43").Value = -4692.2856
# Row 58
$ws.Range("H58").Value = 21687.256
$ws.Range("I58").Value = 292.22223
$ws.Range("J58").Value = 26271.904
$ws.Range("K58").Value = 876.66669
$ws.Range("L58").Value = 78815.712
$ws.Range("M58").Value = -726.66669
$ws.Range("N58").Value = -79115.712
# Row 100
$ws.Range("H100").Value = 2915.5557
$ws.Range("I100").Value = 2250
$ws.Range("J100").Value = 5245
$ws.Range("K100").Value = 2250
$ws.Range("L100").Value = 5245
$ws.Range("M100").Value = -1709
$ws.Range("N100").Value = -6327
# Row 115
$ws.Range("H115").Value = 1725.2858
$ws.Range("I115").Value = 650.8
$ws.Range("J115").Value = 2322.2222
$ws.Range("K115").Value = 1952.4
$ws.Range("L115").Value = 6966.6666
$ws.Range("M115").Value = -385.3999999999999
$ws.Range("N115").Value = -10100.6666
# Row 137
$ws.Range("H137").Value = 17420.262
$ws.Range("I137").Value = 2026.54
$ws.Range("J137").Value = 68732.664
$ws.Range("K137").Value = 6079.62
$ws.Range("L137").Value = 206197.992
$ws.Range("M137").Value = -3529.62
$ws.Range("N137").Value = -211297.992

$ws = $wb.Worksheets.Item("ARM")
# Row 64
$ws.Range("H64").Value = 27416.666
$ws.Range("J64").Value = 27416.666
$ws.Range("L64").Value = 27416.666
$ws.Range("N64").Value = -27912.666
# Row 67
$ws.Range("H67").Value = 27416.666
$ws.Range("J67").Value = 27416.666
$ws.Range("L67").Value = 27416.666
$ws.Range("N67").Value = -29132.666
# Row 74
$ws.Range("H74").Value = 836.6923
$ws.Range("I74").Value = 730.7778
$ws.Range("J74").Value = 1075
$ws.Range("K74").Value = 730.7778
$ws.Range("L74").Value = 1075
$ws.Range("M74").Value = 143.2222
$ws.Range("N74").Value = -2823
# Row 77
$ws.Range("H77").Value = 836.6923
$ws.Range("I77").Value = 730.7778
$ws.Range("J77").Value = 1075
$ws.Range("K77").Value = 3653.889
$ws.Range("L77").Value = 5375
$ws.Range("M77").Value = 714.1110000000003
$ws.Range("N77").Value = -14111
# Row 102
$ws.Range("H102").Value = 2712
$ws.Range("I102").Value = 1801.3334
$ws.Range("J102").Value = 4533.3335
$ws.Range("K102").Value = 1801.3334
$ws.Range("L102").Value = 4533.3335
$ws.Range("M102").Value = -179.3334
$ws.Range("N102").Value = -7777.3335
# Row 124
$ws.Range("H124").Value = 20464.5
$ws.Range("J124").Value = 20464.5
$ws.Range("L124").Value = 20464.5
$ws.Range("N124").Value = -30284.5
# Row 125
$ws.Range("H125").Value = 26081.715
$ws.Range("J125").Value = 26081.715
$ws.Range("L125").Value = 26081.715
$ws.Range("N125").Value = -35921.715
# Row 132
$ws.Range("H132").Value = 22224380
$ws.Range("I132").Value = 40001570
$ws.Range("K132").Value = 120004710
$ws.Range("M132").Value = -120002180

$ws = $wb.Worksheets.Item("BSM")
# Row 18
$ws.Range("H18").Value = 70011
$ws.Range("J18").Value = 70011
$ws.Range("L18").Value = 70011
$ws.Range("N18").Value = -71069
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 133
$ws.Range("H133").Value = 20446.666
$ws.Range("J133").Value = 20446.666
$ws.Range("L133").Value = 20446.666
$ws.Range("N133").Value = -30566.666

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 230209.8
$ws.Range("I31").Value = 1979.2609
$ws.Range("J31").Value = 480176.56
$ws.Range("K31").Value = 1979.2609
$ws.Range("L31").Value = 480176.56
$ws.Range("M31").Value = -1684.2609
$ws.Range("N31").Value = -480766.56
# Row 34
$ws.Range("H34").Value = 230209.8
$ws.Range("I34").Value = 1979.2609
$ws.Range("J34").Value = 480176.56
$ws.Range("K34").Value = 1979.2609
$ws.Range("L34").Value = 480176.56
$ws.Range("M34").Value = -1777.2609
$ws.Range("N34").Value = -480580.56
# Row 99
$ws.Range("H99").Value = 3042.7144
$ws.Range("I99").Value = 1600
$ws.Range("J99").Value = 3283.1667
$ws.Range("K99").Value = 1600
$ws.Range("L99").Value = 3283.1667
$ws.Range("M99").Value = -102
$ws.Range("N99").Value = -6279.1667
# Row 122
$ws.Range("H122").Value = 3484
$ws.Range("I122").Value = 3199.8
$ws.Range("J122").Value = 3578.7334
$ws.Range("K122").Value = 9599.400000000001
$ws.Range("L122").Value = 10736.2002
$ws.Range("M122").Value = -7149.400000000001
$ws.Range("N122").Value = -15636.2002
# Row 126
$ws.Range("H126").Value = 3042.7144
$ws.Range("I126").Value = 1600
$ws.Range("J126").Value = 3283.1667
$ws.Range("K126").Value = 4800
$ws.Range("L126").Value = 9849.500100000001
$ws.Range("M126").Value = -2330
$ws.Range("N126").Value = -14789.5001
# Row 132
$ws.Range("H132").Value = 2159.7354
$ws.Range("I132").Value = 1576.0741
$ws.Range("J132").Value = 4411
$ws.Range("K132").Value = 4728.2223
$ws.Range("L132").Value = 13233
$ws.Range("M132").Value = -2198.2223
$ws.Range("N132").Value = -18293
# Row 134
$ws.Range("H134").Value = 1785.4615
$ws.Range("I134").Value = 800.95654
$ws.Range("J134").Value = 9333.333000000001
$ws.Range("K134").Value = 2402.86962
$ws.Range("L134").Value = 27999.999
$ws.Range("M134").Value = 132.1303800000001
$ws.Range("N134").Value = -33069.999

$ws = $wb.Worksheets.Item("CUL")
# Row 95
$ws.Range("H95").Value = 170004
$ws.Range("I95").Value = 1000024
$ws.Range("J95").Value = 4000
$ws.Range("K95").Value = 3000072
$ws.Range("L95").Value = 12000
$ws.Range("M95").Value = -2998013
$ws.Range("N95").Value = -16118
# Row 114
$ws.Range("H114").Value = 858.65216
$ws.Range("I114").Value = 269
$ws.Range("J114").Value = 982.7895
$ws.Range("K114").Value = 807
$ws.Range("L114").Value = 2948.3685
$ws.Range("M114").Value = 2447
$ws.Range("N114").Value = -9456.3685
# Row 120
$ws.Range("H120").Value = 13353.333
$ws.Range("J120").Value = 20000
$ws.Range("L120").Value = 60000
$ws.Range("N120").Value = -69676
# Row 122
$ws.Range("H122").Value = 1112
$ws.Range("I122").Value = 486
$ws.Range("J122").Value = 1827.4286
$ws.Range("K122").Value = 4374
$ws.Range("L122").Value = 16446.8574
$ws.Range("M122").Value = -1924
$ws.Range("N122").Value = -21346.8574
# Row 130
$ws.Range("H130").Value = 1987.5
$ws.Range("I130").Value = 950
$ws.Range("J130").Value = 2333.3333
$ws.Range("K130").Value = 2850
$ws.Range("L130").Value = 6999.999899999999
$ws.Range("M130").Value = 2170
$ws.Range("N130").Value = -17039.9999
# Row 131
$ws.Range("H131").Value = 2017.3077
$ws.Range("I131").Value = 5770
$ws.Range("J131").Value = 1335
$ws.Range("K131").Value = 17310
$ws.Range("L131").Value = 4005
$ws.Range("M131").Value = -12270
$ws.Range("N131").Value = -14085
# Row 138
$ws.Range("H138").Value = 3787.889
$ws.Range("I138").Value = 1388.2
$ws.Range("J138").Value = 6787.5
$ws.Range("K138").Value = 4164.6
$ws.Range("L138").Value = 20362.5
$ws.Range("M138").Value = 975.3999999999996
$ws.Range("N138").Value = -30642.5

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2070.1875
$ws.Range("I132").Value = 1635.4783
$ws.Range("J132").Value = 2470.12
$ws.Range("K132").Value = 4906.4349
$ws.Range("L132").Value = 7410.36
$ws.Range("M132").Value = -2376.4349
$ws.Range("N132").Value = -12470.36

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1710.6428
$ws.Range("I46").Value = 399.83334
$ws.Range("J46").Value = 2693.75
$ws.Range("K46").Value = 399.83334
$ws.Range("L46").Value = 2693.75
$ws.Range("M46").Value = -211.83334
$ws.Range("N46").Value = -3069.75
# Row 132
$ws.Range("H132").Value = 2572
$ws.Range("I132").Value = 1761.2632
$ws.Range("K132").Value = 5283.7896
$ws.Range("M132").Value = -2753.7896

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 324570.72
$ws.Range("I122").Value = 401696.75
$ws.Range("J122").Value = 3212.1667
$ws.Range("K122").Value = 1205090.25
$ws.Range("L122").Value = 9636.500100000001
$ws.Range("M122").Value = -1202640.25
$ws.Range("N122").Value = -14536.5001
# Row 132
$ws.Range("H132").Value = 10689.1
$ws.Range("I132").Value = 2403.4883
$ws.Range("J132").Value = 31646.824
$ws.Range("K132").Value = 7210.4649
$ws.Range("L132").Value = 94940.47200000001
$ws.Range("M132").Value = -4680.4649
$ws.Range("N132").Value = -100000.472
# Row 136
$ws.Range("H136").Value = 1175.4634
$ws.Range("I136").Value = 787.5517
$ws.Range("J136").Value = 2112.9167
$ws.Range("K136").Value = 2362.6551
$ws.Range("L136").Value = 6338.750100000001
$ws.Range("M136").Value = 187.3449000000001
$ws.Range("N136").Value = -11438.7501
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
# Row 141
$ws.Range("H141").Value = 28200
$ws.Range("J141").Value = 28200
$ws.Range("L141").Value = 28200
$ws.Range("N141").Value = -38560
